$d = $word.ActiveDocument

# The document ends with a "To Do" heading paragraph followed by a single
# trailing empty paragraph (style "para2"). We insert four new "para2"
# paragraphs of text right before that trailing empty paragraph, leaving
# the trailing paragraph itself untouched.
$lastIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($lastIndex)
$r = $lastPara.Range
$r.Collapse(1)

$newText = "dumpProperties(<symbol>) cf. <use>`r" `
    + "Defer widget construction (setTimeout) in case <use> is slow to load`r" `
    + "Try to create problem when removing class=class. (Could look in c-t git history for introduction of that line, and discord chat about the problem.)`r" `
    + ".md`r"

$r.InsertBefore($newText)
